# Update the "Team of Outs" player table on Sheet1.
# The roster/positions/teams are replaced with a new list (15 rows instead
# of 16), so the old data region (A2:C17) is cleared and rewritten, then
# the now-empty trailing row is removed so the sheet dimension shrinks
# from A1:C17 to A1:C16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$players = @(
    "Ja Morant",
    "Isaiah Collier",
    "De'Aaron Fox",
    "Luka Doncic",
    "Scottie Barnes",
    "Mikal Bridges",
    "Miles Bridges",
    "Evan Mobley",
    "Nikola Vucevic",
    "Brook Lopez",
    "Josh Giddey",
    "DeMar DeRozan",
    "Tyler Herro",
    "Royce O'Neale",
    "P.J. Washington"
)

$positions = @(
    "PG",
    "PG,SG",
    "PG",
    "PG,SG",
    "PG,SG,SF,PF",
    "SG,SF,PF",
    "SF,PF",
    "PF,C",
    "PF,C",
    "C",
    "PG,SG,SF",
    "SF,PF",
    "PG,SG",
    "SF,PF",
    "SF,PF"
)

$teams = @(
    "Memphis Grizzlies",
    "Utah Jazz",
    "San Antonio Spurs",
    "Los Angeles Lakers",
    "Toronto Raptors",
    "New York Knicks",
    "Charlotte Hornets",
    "Cleveland Cavaliers",
    "Chicago Bulls",
    "Milwaukee Bucks",
    "Chicago Bulls",
    "Sacramento Kings",
    "Miami Heat",
    "Phoenix Suns",
    "Dallas Mavericks"
)

# Clear out the previous 16 rows of data (rows 2-17).
$ws.Range("A2:C17").ClearContents()

# Write the new data, one row at a time, starting at row 2.
for ($i = 0; $i -lt $players.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $players[$i]
    $ws.Cells.Item($r, 2).Value = $positions[$i]
    $ws.Cells.Item($r, 3).Value = $teams[$i]
}

# The new table only has 15 data rows (2-16); remove the now-blank row 17
# so the sheet dimension becomes A1:C16.
$ws.Rows(17).Delete()
